# Update points 09876543 -> 0.00
#
# The sheet has a row whose "phone" column holds the text value
# "09876543" (leading zero, stored as text). This edit inserts a new
# row just above it that holds the corrected/normalized numeric phone
# number 9876543 with its points (total_points) reset to 0, while the
# original "09876543" text row is preserved unchanged, shifted down by
# one row.
#
# Net effect vs. the original sheet (rows 1-41 untouched):
#   row 42 (new)   -> A=9876543 (number), B=<blank>, C=0
#   row 43 (was 42)-> A="09876543" (text), B=<blank>, C=0
#   dimension grows from A1:C42 to A1:C43

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the "09876543" row (row 42), pushing the
# existing row 42 (and everything below it) down to row 43.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the corrected numeric phone
# number and zeroed points; the birthday column (B) stays blank, same
# as the row it was copied from.
$ws.Range("A42").Value = 9876543
$ws.Range("B42").Value = ""
$ws.Range("C42").Value = 0
